# Commit: "add the NA's under duplicate_image_filename"
#
# Column E (header "duplicate_image_filename" in E1) is blank for the
# practice/generic/unique stimuli rows (rows 2-21). Fill each of those
# rows with the value "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
